$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "mark_red" column (L) with 1 for the newly marked_red sessions
# (rows 26-45), preserving each cell's existing number format/style.
$ws.Range("L26:L45").Value = 1

# Reflect the author's new view/selection state: the window was scrolled
# down (top-left cell now A13) and the active cell/selection moved to A46.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A46").Select()
